# Results from R script
# Updates the last row of price history (row 180) and appends two newer
# trading rows (181, 182) to Sheet 1 of the ESF.MI workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: write a text value into a cell so it ends up stored as a
# shared string (t="s"), matching column G which always holds the
# "close" price repeated as text (rather than being auto-coerced to a
# number, which is what a plain numeric-looking string would become).
# A text formula is entered, then copy/paste-special-values collapses
# it down to a literal shared-string cell with no leftover formula and
# no leftover style/number-format override.
# ---------------------------------------------------------------------
function Set-TextValue($cell, $text) {
    $cell.Formula = "=""" + $text + """"
    $cell.Copy()
    $cell.PasteSpecial(-4163)   # xlPasteValues
    $excel.CutCopyMode = $false
}

# --- Row 180: overwrite the existing last row with the corrected values ---
$ws.Range("A180").Value = 45454.2916666667
$ws.Range("B180").Value = 1200
$ws.Range("C180").Value = 4.44000005722046
$ws.Range("F180").Value = 4.26000022888184
Set-TextValue $ws.Range("G180") "4.26000022888184"

# --- Row 181: new row inserted after the (updated) row 180 ---
$ws.Range("A179").Copy($ws.Range("A181"))
$ws.Range("A181").Value = 45455.2916666667
$ws.Range("B181").Value = 102
$ws.Range("C181").Value = 4.32999992370605
$ws.Range("D181").Value = 4.17000007629395
$ws.Range("E181").Value = 4.26000022888184
$ws.Range("F181").Value = 4.17000007629395
Set-TextValue $ws.Range("G181") "4.17000007629395"
$ws.Range("H181").Value = "ESF.MI"

# --- Row 182: newest row appended at the end ---
$ws.Range("A179").Copy($ws.Range("A182"))
$ws.Range("A182").Value = 45456.6178587963
$ws.Range("B182").Value = 6065
$ws.Range("C182").Value = 4.5
$ws.Range("D182").Value = 4.21999979019165
$ws.Range("E182").Value = 4.44000005722046
$ws.Range("F182").Value = 4.34000015258789
Set-TextValue $ws.Range("G182") "4.34000015258789"
$ws.Range("H182").Value = "ESF.MI"
